# Auto-generated script applying scheduled market-price refresh updates
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 67220.336
$ws.Range("I2").Value = 185
$ws.Range("J2").Value = 167773.33
$ws.Range("K2").Value = 185
$ws.Range("L2").Value = 167773.33
$ws.Range("M2").Value = -72
$ws.Range("N2").Value = -167999.33
$ws.Range("H17").Value = 1197
$ws.Range("J17").Value = 1197
$ws.Range("L17").Value = 3591
$ws.Range("N17").Value = -3927
$ws.Range("H18").Value = 19542.715
$ws.Range("I18").Value = 22559.8
$ws.Range("K18").Value = 22559.8
$ws.Range("M18").Value = -22275.8
$ws.Range("H29").Value = 3785.4285
$ws.Range("I29").Value = 1749.75
$ws.Range("J29").Value = 6499.6665
$ws.Range("K29").Value = 5249.25
$ws.Range("L29").Value = 19498.9995
$ws.Range("M29").Value = -4968.25
$ws.Range("N29").Value = -20060.9995
$ws.Range("H100").Value = 3423.9443
$ws.Range("I100").Value = 3279.6155
$ws.Range("J100").Value = 3799.2
$ws.Range("K100").Value = 3279.6155
$ws.Range("L100").Value = 3799.2
$ws.Range("M100").Value = -2738.6155
$ws.Range("N100").Value = -4881.2
$ws.Range("H132").Value = 1904.1765
$ws.Range("I132").Value = 1872.5
$ws.Range("J132").Value = 2103.2856
$ws.Range("K132").Value = 5617.5
$ws.Range("L132").Value = 6309.8568
$ws.Range("M132").Value = -3087.5
$ws.Range("N132").Value = -11369.8568
$ws.Range("H137").Value = 71999.31
$ws.Range("I137").Value = 1524.4459
$ws.Range("K137").Value = 4573.3377
$ws.Range("M137").Value = -2023.3377
$ws.Range("H138").Value = 1928.2911
$ws.Range("I138").Value = 1348.3334
$ws.Range("K138").Value = 4045.0002
$ws.Range("M138").Value = 1094.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1729
$ws.Range("I25").Value = 1729
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1729
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1327
$ws.Range("N25").ClearContents()
$ws.Range("H74").Value = 242586
$ws.Range("I74").Value = 261197.97
$ws.Range("J74").Value = 130914.14
$ws.Range("K74").Value = 261197.97
$ws.Range("L74").Value = 130914.14
$ws.Range("M74").Value = -260323.97
$ws.Range("N74").Value = -132662.14
$ws.Range("H77").Value = 242586
$ws.Range("I77").Value = 261197.97
$ws.Range("J77").Value = 130914.14
$ws.Range("K77").Value = 1305989.85
$ws.Range("L77").Value = 654570.7
$ws.Range("M77").Value = -1301621.85
$ws.Range("N77").Value = -663306.7
$ws.Range("H95").Value = 29034.5
$ws.Range("J95").Value = 29034.5
$ws.Range("L95").Value = 29034.5
$ws.Range("N95").Value = -34526.5
$ws.Range("H97").Value = 835.17645
$ws.Range("I97").Value = 452.32
$ws.Range("J97").Value = 1898.6666
$ws.Range("K97").Value = 452.32
$ws.Range("L97").Value = 1898.6666
$ws.Range("M97").Value = 43.68000000000001
$ws.Range("N97").Value = -2890.6666
$ws.Range("H103").Value = 29249.25
$ws.Range("J103").Value = 29249.25
$ws.Range("L103").Value = 29249.25
$ws.Range("N103").Value = -31593.25
$ws.Range("H122").Value = 32844.234
$ws.Range("I122").Value = 2678.0417
$ws.Range("K122").Value = 8034.125100000001
$ws.Range("M122").Value = -5584.125100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 16981454
$ws.Range("I86").Value = 26341182
$ws.Range("J86").Value = 44800.715
$ws.Range("K86").Value = 26341182
$ws.Range("L86").Value = 44800.715
$ws.Range("M86").Value = -26340059
$ws.Range("N86").Value = -47046.715
$ws.Range("H89").Value = 16981454
$ws.Range("I89").Value = 26341182
$ws.Range("J89").Value = 44800.715
$ws.Range("K89").Value = 131705910
$ws.Range("L89").Value = 224003.575
$ws.Range("M89").Value = -131700294
$ws.Range("N89").Value = -235235.575
$ws.Range("H94").Value = 33871.555
$ws.Range("I94").Value = 497.5
$ws.Range("J94").Value = 180717.4
$ws.Range("K94").Value = 497.5
$ws.Range("L94").Value = 180717.4
$ws.Range("M94").Value = -46.5
$ws.Range("N94").Value = -181619.4
$ws.Range("H99").Value = 1382.0385
$ws.Range("I99").Value = 1403.6316
$ws.Range("J99").Value = 1323.4286
$ws.Range("K99").Value = 1403.6316
$ws.Range("L99").Value = 1323.4286
$ws.Range("M99").Value = 94.36840000000007
$ws.Range("N99").Value = -4319.4286
$ws.Range("H107").Value = 11225.595
$ws.Range("I107").Value = 14310.926
$ws.Range("J107").Value = 2895.2
$ws.Range("K107").Value = 14310.926
$ws.Range("L107").Value = 2895.2
$ws.Range("M107").Value = -12390.926
$ws.Range("N107").Value = -6735.2
$ws.Range("H134").Value = 3268.325
$ws.Range("I134").Value = 2865.516
$ws.Range("J134").Value = 4655.778
$ws.Range("K134").Value = 8596.548000000001
$ws.Range("L134").Value = 13967.334
$ws.Range("M134").Value = -6061.548000000001
$ws.Range("N134").Value = -19037.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2155.6047
$ws.Range("I31").Value = 1607.3024
$ws.Range("K31").Value = 1607.3024
$ws.Range("M31").Value = -1312.3024
$ws.Range("H34").Value = 2155.6047
$ws.Range("I34").Value = 1607.3024
$ws.Range("K34").Value = 1607.3024
$ws.Range("M34").Value = -1405.3024
$ws.Range("H58").Value = 3519.2563
$ws.Range("I58").Value = 3078.138
$ws.Range("J58").Value = 4798.5
$ws.Range("K58").Value = 3078.138
$ws.Range("L58").Value = 4798.5
$ws.Range("M58").Value = -2875.138
$ws.Range("N58").Value = -5204.5
$ws.Range("H134").Value = 1866.5938
$ws.Range("I134").Value = 1583.0741
$ws.Range("K134").Value = 4749.2223
$ws.Range("M134").Value = -2214.2223
$ws.Range("H136").Value = 3519.2563
$ws.Range("I136").Value = 3078.138
$ws.Range("J136").Value = 4798.5
$ws.Range("K136").Value = 9234.414000000001
$ws.Range("L136").Value = 14395.5
$ws.Range("M136").Value = -6684.414000000001
$ws.Range("N136").Value = -19495.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2999.5
$ws.Range("J32").Value = 4500
$ws.Range("L32").Value = 13500
$ws.Range("N32").Value = -14066
$ws.Range("H46").Value = 1781.0667
$ws.Range("I46").Value = 871.25
$ws.Range("K46").Value = 2613.75
$ws.Range("M46").Value = -2522.75
$ws.Range("H113").Value = 1487.0526
$ws.Range("I113").Value = 613.5
$ws.Range("J113").Value = 1720
$ws.Range("K113").Value = 1840.5
$ws.Range("L113").Value = 5160
$ws.Range("M113").Value = 329.5
$ws.Range("N113").Value = -9500
$ws.Range("H131").Value = 2958.0967
$ws.Range("J131").Value = 3739.15
$ws.Range("L131").Value = 11217.45
$ws.Range("N131").Value = -21297.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 21297948
$ws.Range("I97").Value = 35748430
$ws.Range("K97").Value = 35748430
$ws.Range("M97").Value = -35747934

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 4967
$ws.Range("J9").Value = 3960.4
$ws.Range("L9").Value = 3960.4
$ws.Range("N9").Value = -4408.4
$ws.Range("H46").Value = 2400.5469
$ws.Range("I46").Value = 797
$ws.Range("K46").Value = 797
$ws.Range("M46").Value = -609
$ws.Range("H68").Value = 3649.5
$ws.Range("I68").Value = 1749.8334
$ws.Range("K68").Value = 1749.8334
$ws.Range("M68").Value = -1000.8334
$ws.Range("H71").Value = 3649.5
$ws.Range("I71").Value = 1749.8334
$ws.Range("K71").Value = 8749.166999999999
$ws.Range("M71").Value = -5005.166999999999
$ws.Range("H93").Value = 1153.6666
$ws.Range("I93").Value = 1252.6111
$ws.Range("J93").Value = 560
$ws.Range("K93").Value = 1252.6111
$ws.Range("L93").Value = 560
$ws.Range("M93").Value = -4.611100000000079
$ws.Range("N93").Value = -3056
$ws.Range("H122").Value = 3651.35
$ws.Range("I122").Value = 3451.75
$ws.Range("K122").Value = 10355.25
$ws.Range("M122").Value = -7905.25
$ws.Range("H136").Value = 50579.844
$ws.Range("I136").Value = 2606.4614
$ws.Range("J136").Value = 154522.17
$ws.Range("K136").Value = 7819.3842
$ws.Range("L136").Value = 463566.51
$ws.Range("M136").Value = -5269.3842
$ws.Range("N136").Value = -468666.51

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8297
$ws.Range("I14").Value = 4390
$ws.Range("K14").Value = 4390
$ws.Range("M14").Value = -4222
$ws.Range("H19").Value = 6699.6665
$ws.Range("I19").Value = 6900
$ws.Range("J19").Value = 6499.3335
$ws.Range("K19").Value = 6900
$ws.Range("L19").Value = 6499.3335
$ws.Range("M19").Value = -6726
$ws.Range("N19").Value = -6847.3335
$ws.Range("H39").Value = 38995
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 38995
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 38995
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -39821
$ws.Range("H100").Value = 52632324
$ws.Range("H122").Value = 2478.3625
$ws.Range("J122").Value = 3284.5334
$ws.Range("L122").Value = 9853.600199999999
$ws.Range("N122").Value = -14753.6002
